# "Handles float input without breaking stuff"
# Update the quiz marksheet: fix the summary numbers (rows 10-12), convert the
# "-1" text in C11 to a real number, fill in the student's actual answers
# (column A, and D17/D18) with correct/incorrect colouring, and drop the
# now-unused "Student Ans / Correct Ans" sub-tables that are no longer part
# of the report (columns D/E for most rows, and the whole G/H block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# --- Row 10: summary header counts -----------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 18
$ws.Range("E10").Value = 28

# --- Row 11: marking scheme -------------------------------------------
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# --- Row 12: totals -----------------------------------------------------
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 28
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "25/112"

# --- Drop the 3rd ("Student Ans" / "Correct Ans") sub-table entirely ----
$ws.Range("G15:H40").Clear()

# --- Row 16: drop the 2nd sub-table entry for this row -------------------
$ws.Range("D16:E16").Clear()

# --- Row 17: student answered Option A (wrong, correct was Option C) -----
$ws.Range("D17").Value = "Option A"
$ws.Range("D17").Style = "incorrectStyle"

# --- Row 18: student filled in Option B for both the 1st and 2nd table ---
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("D18").Value = "Option B"
$ws.Range("D18").Style = "incorrectStyle"

# --- Rows 19-40: drop the 2nd sub-table (D/E) for every remaining row ----
$ws.Range("D19:E40").Clear()

# --- Column A: fill in the rest of the student's attempted answers -------
$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"

$ws.Range("A26").Value = "Option C"
$ws.Range("A26").Style = "correctStyle"

$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"

$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"

$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"

$ws.Range("A34").Value = "Option D"
$ws.Range("A34").Style = "incorrectStyle"

$ws.Range("A37").Value = "Option A"
$ws.Range("A37").Style = "correctStyle"
